# Generate Report for Handoff
# Updates the localization status from "In Translation" to
# "Ready for handoff" and refreshes the associated timestamps across the
# Overview / zh-cn / de-de sheets, widening the Status-related columns to
# fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-09-07 11:22:56"

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-09-07 11:22:51"

# --- de-de sheet ------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-09-07 11:22:56"

# --- Widen the Status columns to fit the longer text -------------------
# ColumnWidth is quantized to the host's pixel grid, so feed it an input
# that lands on the grid point closest to the authored 17.2159881591797
# character-width target (rounds to 17.166666666666668 on this host).
$newColWidth = 16.3
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth
$dede.Columns.Item(3).ColumnWidth = $newColWidth
